$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 28367
$ws.Range("D2").Value = 321
$ws.Range("E2").Value = 6042
$ws.Range("F2").Value = 78
$ws.Range("G2").Value = 24.84
$ws.Range("H2").Value = 26.09
$ws.Range("K2").Value = 24327
$ws.Range("L2").Value = 299

# Row 3
$ws.Range("B3").Value = 44024
$ws.Range("C3").Value = 111597
$ws.Range("D3").Value = 8325
$ws.Range("E3").Value = 10493
$ws.Range("G3").Value = 9.4
$ws.Range("H3").Value = 8.22

# Row 5
$ws.Range("B5").ClearContents()
$ws.Range("B5").ClearFormats()
$ws.Range("C5").ClearContents()
$ws.Range("C5").ClearFormats()
$ws.Range("D5").ClearContents()
$ws.Range("D5").ClearFormats()
$ws.Range("E5").ClearContents()
$ws.Range("E5").ClearFormats()
$ws.Range("F5").ClearContents()
$ws.Range("F5").ClearFormats()
$ws.Range("G5").ClearContents()
$ws.Range("G5").ClearFormats()
$ws.Range("H5").ClearContents()
$ws.Range("H5").ClearFormats()
$ws.Range("I5").Value = $False
$ws.Range("J5").Value = $False
$ws.Range("O5").Value = "An error occurred. ... AssertionError('No percentage deaths found.')"

# Row 6
$ws.Range("B6").Value = 44024
$ws.Range("C6").Value = 153916
$ws.Range("D6").Value = 7187
$ws.Range("E6").Value = 25828
$ws.Range("F6").Value = 1990
$ws.Range("G6").Value = 16.78
$ws.Range("H6").Value = 27.69

# Row 8
$ws.Range("B8").Value = 44024
$ws.Range("C8").Value = 19648
$ws.Range("D8").Value = 184

# Row 10
$ws.Range("B10").Value = 44024
$ws.Range("C10").Value = 19389
$ws.Range("D10").Value = 625
$ws.Range("E10").Value = 2072
$ws.Range("G10").Value = 15.91
$ws.Range("H10").Value = 4.12
$ws.Range("K10").Value = 13029
$ws.Range("L10").Value = 587

# Row 12
$ws.Range("B12").Value = 44024
$ws.Range("C12").Value = 36448
$ws.Range("D12").Value = 820
$ws.Range("E12").Value = 6076
$ws.Range("G12").Value = 18.5
$ws.Range("H12").Value = 23.76
$ws.Range("K12").Value = 32836
$ws.Range("L12").Value = 808

# Row 14
$ws.Range("B14").Value = 44024
$ws.Range("C14").Value = 61960
$ws.Range("D14").Value = 741
$ws.Range("E14").Value = 12551
$ws.Range("F14").Value = 263
$ws.Range("G14").Value = 20.26
$ws.Range("H14").Value = 35.49

# Row 16
$ws.Range("B16").Value = 44024
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "29484"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "215"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "747"

# Row 18
$ws.Range("B18").ClearContents()
$ws.Range("B18").ClearFormats()
$ws.Range("C18").ClearContents()
$ws.Range("C18").ClearFormats()
$ws.Range("D18").ClearContents()
$ws.Range("D18").ClearFormats()
$ws.Range("E18").ClearContents()
$ws.Range("E18").ClearFormats()
$ws.Range("F18").ClearContents()
$ws.Range("F18").ClearFormats()
$ws.Range("G18").ClearContents()
$ws.Range("G18").ClearFormats()
$ws.Range("H18").ClearContents()
$ws.Range("H18").ClearFormats()
$ws.Range("J18").Value = $False
$ws.Range("K18").ClearContents()
$ws.Range("K18").ClearFormats()
$ws.Range("L18").ClearContents()
$ws.Range("L18").ClearFormats()
$ws.Range("O18").Value = "An error occurred. ... AssertionError('Unable to find ArcGIS ID 554ada3bc8b147abad21ae23d4a7ba3a')"

# Row 19
$ws.Range("B19").Value = 44023
$ws.Range("B19").NumberFormat = "YYYY-MM-DD"
$ws.Range("C19").Value = 320804
$ws.Range("D19").Value = 6989
$ws.Range("E19").Value = 9021
$ws.Range("F19").Value = 613
$ws.Range("G19").Value = 4.4
$ws.Range("H19").Value = 8.9
$ws.Range("J19").Value = $True
$ws.Range("K19").Value = 206109
$ws.Range("L19").Value = 6888
$ws.Range("O19").Value = "Success!"

# Row 20
$ws.Range("B20").Value = 44024
$ws.Range("C20").Value = 15028
$ws.Range("D20").Value = 545
$ws.Range("E20").Value = 282
$ws.Range("G20").Value = 1.88

# Row 21
$ws.Range("B21").Value = 44024
$ws.Range("C21").Value = 1479
$ws.Range("E21").Value = 35
$ws.Range("G21").Value = 1.4
$ws.Range("K21").Value = 2499

# Row 22
$ws.Range("B22").Value = 44024
$ws.Range("C22").Value = 25438
$ws.Range("D22").Value = 719
$ws.Range("E22").Value = 2914
$ws.Range("G22").Value = 0.11
$ws.Range("H22").Value = 0.13

# Row 24
$ws.Range("B24").Value = 44024
$ws.Range("C24").Value = 19929
$ws.Range("E24").Value = 736
$ws.Range("G24").Value = 4.65
$ws.Range("K24").Value = 15828

# Row 28
$ws.Range("B28").Value = 44023
$ws.Range("C28").Value = 133549
$ws.Range("D28").Value = 3809
$ws.Range("E28").Value = 3554
$ws.Range("F28").Value = 387
$ws.Range("G28").Value = 4.7
$ws.Range("H28").Value = 10.92
$ws.Range("K28").Value = 75586
$ws.Range("L28").Value = 3543

# Row 33
$ws.Range("B33").Value = 44024
$ws.Range("C33").Value = 21172
$ws.Range("E33").Value = 1252
$ws.Range("G33").Value = 7.62
$ws.Range("H33").Value = 8.06
$ws.Range("K33").Value = 16440
$ws.Range("L33").Value = 273

# Row 34
$ws.Range("B34").Value = 44024
$ws.Range("C34").Value = 116926
$ws.Range("D34").Value = 3001
$ws.Range("E34").Value = 31278
$ws.Range("F34").Value = 1398
$ws.Range("G34").Value = 26.75
$ws.Range("H34").Value = 46.58

# Row 36
$ws.Range("B36").Value = 44024
$ws.Range("C36").Value = 10902
$ws.Range("E36").Value = 150
$ws.Range("G36").Value = 1.38

# Row 38
$ws.Range("B38").Value = 44024
$ws.Range("C38").Value = 36913
$ws.Range("E38").Value = 1890
$ws.Range("G38").Value = 6.33
$ws.Range("K38").Value = 29866
$ws.Range("L38").Value = 1661

# Row 41
$ws.Range("B41").Value = 44024
$ws.Range("C41").Value = 69250
$ws.Range("D41").Value = 5984
$ws.Range("E41").Value = 20548
$ws.Range("G41").Value = 29.67
$ws.Range("H41").Value = 39.94
